$wb = $excel.ActiveWorkbook

# --- Add the new "StudyChairs" worksheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "StudyChairs"

# --- Header row ---
$ws.Range("A1").Value = "Model Name"
$ws.Range("B1").Value = "Brand Name"
$ws.Range("C1").Value = "Price"

# Build the bold + centered header style once on A1, then propagate the
# exact same style to B1/C1 via a format-only paste so no intermediate
# per-cell style records are left behind in the styles table.
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Copy()
$ws.Range("B1:C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Data rows (entered in the same row-major order as the source data) ---
$ws.Range("B2").Value = "By Urban Ladder"
$ws.Range("C2").Value = "₹12,287"

$ws.Range("A3").Value = "Galen Study Chair In Black Colour"
$ws.Range("B3").Value = "By Urban Ladder"
$ws.Range("C3").Value = "₹7,505"

$ws.Range("A4").Value = "Hawley Study Chair"
$ws.Range("B4").Value = "By Urban Ladder"
$ws.Range("C4").Value = "₹6,440"

# A2 has no model name recorded -- represent it as an (empty) text cell
# rather than leaving it completely blank.
$ws.Range("A2").Value = "'"
$ws.Range("A2").Style = "Normal"

# --- Column widths ---
$ws.Columns.Item(1).ColumnWidth = 25.5
$ws.Columns.Item(2).ColumnWidth = 34.16666666666667
$ws.Columns.Item(3).ColumnWidth = 25.333333333333332

# --- Selection / active cell matches the saved view state ---
$ws.Range("A2").Select()

Write-Host "StudyChairs sheet added"
